$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "Acknowledge (Backchannel)"
$ws.Range("I9").Value = "sd"
$ws.Range("J9").Value = "Statement-non-opinion"
$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"
$ws.Range("I20").Value = "aa"
$ws.Range("J20").Value = "Agree/Accept"
$ws.Range("I22").Value = "ba"
$ws.Range("J22").Value = "Appreciation"
$ws.Range("I24").Value = "ba"
$ws.Range("J24").Value = "Appreciation"
$ws.Range("I28").Value = "sd"
$ws.Range("J28").Value = "Statement-non-opinion"
$ws.Range("I34").Value = "sd"
$ws.Range("J34").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"
$ws.Range("I42").Value = "%"
$ws.Range("J42").Value = "Uninterpretable"
$ws.Range("I45").Value = "sv"
$ws.Range("J45").Value = "Statement-opinion"
$ws.Range("I51").Value = "aa"
$ws.Range("J51").Value = "Agree/Accept"
$ws.Range("I68").Value = "ba"
$ws.Range("J68").Value = "Appreciation"
$ws.Range("I78").Value = "sd"
$ws.Range("J78").Value = "Statement-non-opinion"
$ws.Range("I91").Value = "sd"
$ws.Range("J91").Value = "Statement-non-opinion"
$ws.Range("I105").Value = "aa"
$ws.Range("J105").Value = "Agree/Accept"
$ws.Range("I113").Value = "sd"
$ws.Range("J113").Value = "Statement-non-opinion"
$ws.Range("I116").Value = "b"
$ws.Range("J116").Value = "Acknowledge (Backchannel)"
$ws.Range("I119").Value = "ba"
$ws.Range("J119").Value = "Appreciation"
$ws.Range("I120").Value = "b"
$ws.Range("J120").Value = "Acknowledge (Backchannel)"
$ws.Range("I153").Value = "sv"
$ws.Range("J153").Value = "Statement-opinion"
$ws.Range("I170").Value = "b"
$ws.Range("J170").Value = "Acknowledge (Backchannel)"
$ws.Range("I171").Value = "sd"
$ws.Range("J171").Value = "Statement-non-opinion"
$ws.Range("I172").Value = "sv"
$ws.Range("J172").Value = "Statement-opinion"
$ws.Range("I183").Value = "%"
$ws.Range("J183").Value = "Uninterpretable"
$ws.Range("I185").Value = "sd"
$ws.Range("J185").Value = "Statement-non-opinion"
$ws.Range("I202").Value = "aa"
$ws.Range("J202").Value = "Agree/Accept"
$ws.Range("I203").Value = "sd"
$ws.Range("J203").Value = "Statement-non-opinion"
$ws.Range("I206").Value = "b"
$ws.Range("J206").Value = "Acknowledge (Backchannel)"
$ws.Range("I214").Value = "aa"
$ws.Range("J214").Value = "Agree/Accept"
$ws.Range("I220").Value = "sd"
$ws.Range("J220").Value = "Statement-non-opinion"
$ws.Range("I223").Value = "sv"
$ws.Range("J223").Value = "Statement-opinion"
$ws.Range("I225").Value = "sd"
$ws.Range("J225").Value = "Statement-non-opinion"
$ws.Range("I234").Value = "sv"
$ws.Range("J234").Value = "Statement-opinion"
$ws.Range("I235").Value = "sd"
$ws.Range("J235").Value = "Statement-non-opinion"
$ws.Range("I241").Value = "sd"
$ws.Range("J241").Value = "Statement-non-opinion"
$ws.Range("I267").Value = "sd"
$ws.Range("J267").Value = "Statement-non-opinion"
$ws.Range("I273").Value = "sv"
$ws.Range("J273").Value = "Statement-opinion"
$ws.Range("I282").Value = "sd"
$ws.Range("J282").Value = "Statement-non-opinion"
$ws.Range("I283").Value = "sd"
$ws.Range("J283").Value = "Statement-non-opinion"
$ws.Range("I291").Value = "sv"
$ws.Range("J291").Value = "Statement-opinion"
$ws.Range("I293").Value = "sd"
$ws.Range("J293").Value = "Statement-non-opinion"
$ws.Range("I301").Value = "ba"
$ws.Range("J301").Value = "Appreciation"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I329").Value = "sv"
$ws.Range("J329").Value = "Statement-opinion"
$ws.Range("I337").Value = "qy"
$ws.Range("J337").Value = "Yes-No-Question"
$ws.Range("I341").Value = "b"
$ws.Range("J341").Value = "Acknowledge (Backchannel)"
$ws.Range("I343").Value = "%"
$ws.Range("J343").Value = "Uninterpretable"
$ws.Range("I344").Value = "sv"
$ws.Range("J344").Value = "Statement-opinion"
$ws.Range("I347").Value = "sv"
$ws.Range("J347").Value = "Statement-opinion"
$ws.Range("I352").Value = "sd"
$ws.Range("J352").Value = "Statement-non-opinion"
$ws.Range("I370").Value = "sd"
$ws.Range("J370").Value = "Statement-non-opinion"
$ws.Range("I371").Value = "ba"
$ws.Range("J371").Value = "Appreciation"
$ws.Range("I377").Value = "aa"
$ws.Range("J377").Value = "Agree/Accept"
$ws.Range("I387").Value = "sv"
$ws.Range("J387").Value = "Statement-opinion"
$ws.Range("I389").Value = "b"
$ws.Range("J389").Value = "Acknowledge (Backchannel)"
$ws.Range("I393").Value = "aa"
$ws.Range("J393").Value = "Agree/Accept"
$ws.Range("I394").Value = "ba"
$ws.Range("J394").Value = "Appreciation"
$ws.Range("I401").Value = "sv"
$ws.Range("J401").Value = "Statement-opinion"
$ws.Range("I406").Value = "aa"
$ws.Range("J406").Value = "Agree/Accept"
$ws.Range("I416").Value = "sv"
$ws.Range("J416").Value = "Statement-opinion"
$ws.Range("I423").Value = "sv"
$ws.Range("J423").Value = "Statement-opinion"
$ws.Range("I426").Value = "sv"
$ws.Range("J426").Value = "Statement-opinion"
$ws.Range("I431").Value = "sd"
$ws.Range("J431").Value = "Statement-non-opinion"
$ws.Range("I468").Value = "sv"
$ws.Range("J468").Value = "Statement-opinion"
$ws.Range("I482").Value = "ba"
$ws.Range("J482").Value = "Appreciation"
$ws.Range("I489").Value = "sd"
$ws.Range("J489").Value = "Statement-non-opinion"
$ws.Range("I496").Value = "aa"
$ws.Range("J496").Value = "Agree/Accept"
$ws.Range("I499").Value = "sd"
$ws.Range("J499").Value = "Statement-non-opinion"
$ws.Range("I509").Value = "sv"
$ws.Range("J509").Value = "Statement-opinion"
$ws.Range("I520").Value = "sd"
$ws.Range("J520").Value = "Statement-non-opinion"
$ws.Range("I541").Value = "%"
$ws.Range("J541").Value = "Uninterpretable"
$ws.Range("I544").Value = "b"
$ws.Range("J544").Value = "Acknowledge (Backchannel)"
$ws.Range("I546").Value = "ba"
$ws.Range("J546").Value = "Appreciation"
$ws.Range("I561").Value = "aa"
$ws.Range("J561").Value = "Agree/Accept"
$ws.Range("I563").Value = "aa"
$ws.Range("J563").Value = "Agree/Accept"
$ws.Range("I566").Value = "aa"
$ws.Range("J566").Value = "Agree/Accept"
$ws.Range("I581").Value = "sv"
$ws.Range("J581").Value = "Statement-opinion"
$ws.Range("I583").Value = "aa"
$ws.Range("J583").Value = "Agree/Accept"
$ws.Range("I584").Value = "aa"
$ws.Range("J584").Value = "Agree/Accept"
